$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Row 25: new event row ------------------------------------------------
# Columns B..K and P are left blank (present, empty cells) for this event,
# mirroring how the other rows leave unused measurement columns blank.
$blankRange = $ws.Range("B25:K25")
$blankRange.Style = "Normal"
$ws.Range("P25").Style = "Normal"

# A25 = "24" (card number) -- written through a helper cell + copy/paste
# special so the value lands as text ("24") instead of being auto-coerced
# to a number by the normal Value assignment heuristic.
$helper = $ws.Cells.Item(200, 1)
$helper.Formula = '="" & "24"'
$helper.Copy()
$ws.Range("A25").PasteSpecial(-4163)
$helper.ClearContents()

# L25..O25 hold the actual event data.
$ws.Range("L25").Value = "18/12/2025"
$ws.Range("M25").Value = "كسره جلبه بليه باب كرد امامي"
$ws.Range("N25").Value = "تم تغير جلبه بليه"
$ws.Range("O25").Value = "ابراهيم ،ناجي"

$excel.CutCopyMode = $false
